$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.084.05'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.020.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.03%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.38'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +13.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.015.62'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.69%  '
$ws.Range('E11').Value = '  +5.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.85'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.04%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.004.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.520.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.99%  '
$ws.Range('E18').Value = '  +7.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.013.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '458.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.36%  '
$ws.Range('E22').Value = '  +6.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.43'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.67%  '
$ws.Range('E25').Value = '  +13.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.49'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.66'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.22%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +16.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +17.21%  '
$ws.Range('E31').Value = '  -6.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.62'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.82%  '
$ws.Range('E34').Value = '  +4.89%  '
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.64%  '
$ws.Range('E37').Value = '  +8.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.19'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.73%  '
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('E41').Value = '  +17.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.123'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '395.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0361'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.803.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.86%  '
$ws.Range('E51').Value = '  +4.44%  '
